$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on every cell we touch so Excel does not
# auto-convert numeric-looking or date-looking strings into
# real numbers/dates (the source file stores them as plain text).
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = 'AVINA PRODUCE TIN: xxxxx8949'
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = '546.52'
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = '7888'
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = '2024-11-06'
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = 'AVINA PRODUCE INC'
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '72.92'
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = '2024-11-06'
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = 'AVINA PRODUCE, INC.'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '546.52'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '72.92'
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = 'AVINA PRODUCE TIN: xxxxx8949'
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = '546.52'
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = '7888'
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = 'AVINA PRODUCE INC'
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '72.92'
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = 'AVINA PRODUCE, INC.'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '546.52'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '72.92'
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = 'AVINA PRODUCE TIN: xxxxx8949'
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = '546.52'
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = '7888'
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = 'AVINA PRODUCE INC'
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '72.92'
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = 'AVINA PRODUCE, INC.'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '546.52'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '72.92'
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = 'AVINA PRODUCE TIN: xxxxx8949'
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = '546.52'
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = '7888'
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'AVINA PRODUCE INC'
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '72.92'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'AVINA PRODUCE, INC.'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '546.52'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '72.92'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'AVINA PRODUCE'
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = '1,426.00'
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = '7888'
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = 'AVINA PRODUCE INC'
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '72.92'
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = 'AVINA PRODUCE, INC.'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '546.52'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '72.92'
